$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1: new column headers "I0" and "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting of the existing header cells (bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for the new I (I0) and J (IF) columns, rows 2-53
$iValues = @(9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,8,9,9,9,9,9,10,9,9,9,9,9,8,9,9,9,10,9,9,9,9,9,9,8,8,8,4,6,9)
$jValues = @(9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,10,9,9,8,9,9,9,9,9,10,9,9,9,9,9,8,9,9,9,10,9,9,9,9,9,9,8,8,8,4,6,9)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
